$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts old N/O/P -> O/P/Q)
$ws.Columns("N:N").Insert()

# The newly inserted column inherits the width of the column to its left (M)
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Switch the active sheet/selection to "Repayment schedule" (was "Transactions")
$ws.Activate()
$ws.Range("R7").Select()
